# Update "想去人数" (interested-people count) figures that were refreshed
# by the site's scheduled data pull (gh-pages output regenerated at 456a3b4).
#
# Sheet 1 ("展览") and Sheet 4 ("全部类型") both list the same events, so the
# same six events need their F-column counts bumped on both sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 6599
$ws1.Range("F8").Value = 96
$ws1.Range("F12").Value = 390
$ws1.Range("F15").Value = 3281
$ws1.Range("F18").Value = 1919
$ws1.Range("F19").Value = 34

# Sheet 4: 全部类型 (All types) - same events, rows offset by one extra row
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 6599
$ws4.Range("F9").Value = 96
$ws4.Range("F13").Value = 390
$ws4.Range("F16").Value = 3281
$ws4.Range("F19").Value = 1919
$ws4.Range("F20").Value = 34
